$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Al faro"
$ws.Range("B6").Value = "Virginia Woolf"
$ws.Range("C6").Value = "Lumen"
